# Update the "Estado de Cuenta" worksheet with refreshed database values:
#  - Total "VALOR MORA" (E11) and "Cant. Periodos" (F13) summary figures.
#  - Worker detail rows (16-18): swap out KATTY JULIETH PALACIO PANTOJA's two
#    rows for JESUS MANUEL BARRIOS MURILLO (row 16, moved up) and a brand new
#    worker CESAR MAURICIO BERRIO PEREZ (row 17), and bump row 18's period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Summary block ---
$ws.Range("E11").Value = 186820
$ws.Range("F13").Value = 2

# --- Row 16: CC / 1137220831 / JESUS MANUEL BARRIOS MURILLO / 2507 ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1137220831"
$ws.Range("D16").Value = "JESUS MANUEL BARRIOS MURILLO"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 64940
$ws.Range("G16").Value = 1623500

# --- Row 17: CC / 1002247155 / CESAR MAURICIO BERRIO PEREZ / 2508 ---
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1002247155"
$ws.Range("D17").Value = "CESAR MAURICIO BERRIO PEREZ"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# --- Row 18: CC / 1137220831 / JESUS MANUEL BARRIOS MURILLO / 2508 ---
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1137220831"
$ws.Range("D18").Value = "JESUS MANUEL BARRIOS MURILLO"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 64940
$ws.Range("G18").Value = 1623500
